# Apply the latest crypto price/volume(1h) snapshot to the sheet.
# Values that are purely numeric-looking (e.g. "105.18") are written with a
# leading apostrophe so Excel keeps them as text, matching the source data
# which stores prices/percentages as text (to preserve formats like
# "51.659.97" or "  -0.41%  " verbatim).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.659.97'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '2.941.01'
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '''380.37'
$ws.Range("E5").Value = '  +7.51%  '
$ws.Range("D6").Value = '''105.18'
$ws.Range("E6").Value = '  -1.40%  '
$ws.Range("E7").Value = '  -1.64%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '''0.595'
$ws.Range("E9").Value = '  -1.06%  '
$ws.Range("D10").Value = '''37.25'
$ws.Range("E10").Value = '  -1.13%  '
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("E12").Value = '  -0.97%  '
$ws.Range("D13").Value = '''18.58'
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("D14").Value = '3.398.10'
$ws.Range("E14").Value = '  -1.34%  '
$ws.Range("D15").Value = '''7.45'
$ws.Range("E15").Value = '  -0.66%  '
$ws.Range("D16").Value = '2.936.84'
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D17").Value = '''0.953'
$ws.Range("E17").Value = '  -3.48%  '
$ws.Range("D18").Value = '51.581.06'
$ws.Range("E18").Value = '  -0.47%  '
$ws.Range("E19").Value = '  +1.60%  '
$ws.Range("D20").Value = '''7.40'
$ws.Range("E20").Value = '  +0.49%  '
$ws.Range("D21").Value = '''13.17'
$ws.Range("E21").Value = '  -1.37%  '
$ws.Range("D22").Value = '0.0₃0956'
$ws.Range("E22").Value = '  -0.61%  '
$ws.Range("D23").Value = '''68.53'
$ws.Range("E23").Value = '  -0.64%  '
$ws.Range("D24").Value = '''263.14'
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("E25").Value = '  +4.68%  '
$ws.Range("D26").Value = '''7.42'
$ws.Range("E26").Value = '  +18.89%  '
$ws.Range("E27").Value = '  -4.36%  '
$ws.Range("D28").Value = '''0.169'
$ws.Range("E28").Value = '  -4.46%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").Value = '''7.43'
$ws.Range("E30").Value = '  +2.01%  '
$ws.Range("D31").Value = '''25.93'
$ws.Range("E31").Value = '  -2.85%  '
$ws.Range("E32").Value = '  -7.35%  '
$ws.Range("D33").Value = '''9.89'
$ws.Range("E33").Value = '  -1.74%  '
$ws.Range("D34").Value = '''52.28'
$ws.Range("E34").Value = '  +2.87%  '
$ws.Range("D35").Value = '''34.39'
$ws.Range("E35").Value = '  -3.75%  '
$ws.Range("E36").Value = '  -4.45%  '
$ws.Range("D37").Value = '''0.0435'
$ws.Range("E37").Value = '  +1.80%  '
$ws.Range("E38").Value = '  +0.36%  '
$ws.Range("D39").Value = '''3.04'
$ws.Range("E39").Value = '  -5.51%  '
$ws.Range("D40").Value = '''17.32'
$ws.Range("E40").Value = '  -0.16%  '
$ws.Range("D41").Value = '''2.64'
$ws.Range("E41").Value = '  -6.51%  '
$ws.Range("D42").Value = '''1.84'
$ws.Range("E42").Value = '  -4.54%  '
$ws.Range("E43").Value = '  -0.52%  '
$ws.Range("D44").Value = '''122.93'
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("D45").Value = '''21.89'
$ws.Range("E45").Value = '  -4.39%  '
$ws.Range("E46").Value = '  -3.63%  '
$ws.Range("D47").Value = '''0.278'
$ws.Range("E47").Value = '  +16.97%  '
$ws.Range("D48").Value = '2.027.42'
$ws.Range("E48").Value = '  -3.67%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").Value = '''3.21'
$ws.Range("E50").Value = '  -1.69%  '
$ws.Range("D51").Value = '''0.0324'
$ws.Range("E51").Value = '  -0.73%  '
